$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new "add git repo to existing R project" section below existing content ---

# Section heading (row 26) - same plain formatting as other section headers (e.g. B19, B23)
$ws.Range("B26").Value = "add git repo to existing R project"

# Step detail lines (rows 27-31) - indented like the other step/sub-item rows (e.g. B20:B21, B24)
$ws.Range("B27").Value = "open the R project & go to console"
$ws.Range("B28").Value = "library(usethis)"
$ws.Range("B29").Value = "use_git()"
$ws.Range("B30").Value = "follow the user prompts - confirm to add files to the repo"
$ws.Range("B31").Value = "restart Rstudio & look for git tab in upper right block"

$indentedRange = $ws.Range("B27:B31")
$indentedRange.HorizontalAlignment = -4131  # xlLeft
$indentedRange.IndentLevel = 1

# --- Update the sheet view to reflect where the user ended up after the edit ---
$ws.Range("B33").Select()
$excel.ActiveWindow.ScrollRow = 9
$excel.ActiveWindow.ScrollColumn = 1
